$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "ODI Bowling Extra" worksheet after "ODI Batting Extra"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "ODI Bowling Extra"

# Header row
$ws5.Range("A1").Value = "MATCH_CODE"
$ws5.Range("B1").Value = "MAIDEN_OVERS"
$ws5.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

$hdr5 = $ws5.Range("A1:C1")
$hdr5.Font.Bold = $true
$hdr5.HorizontalAlignment = -4108
$hdr5.VerticalAlignment = -4160
$hdr5.Borders.LineStyle = 1

# Data rows (values stored as text, matching the source sheets)

$ws5.Range("A2").Value = "'4335"
$ws5.Range("B2").Value = "'0"
$ws5.Range("C2").Value = "'"
$ws5.Range("A3").Value = "'4340"
$ws5.Range("B3").Value = "'0"
$ws5.Range("C3").Value = "'10.00%"
$ws5.Range("A4").Value = "'4348"
$ws5.Range("B4").Value = "'0"
$ws5.Range("C4").Value = "'10.00%"
$ws5.Range("A5").Value = "'4377"
$ws5.Range("B5").Value = "'0"
$ws5.Range("C5").Value = "'"
$ws5.Range("A6").Value = "'4378"
$ws5.Range("B6").Value = "'0"
$ws5.Range("C6").Value = "'10.00%"
$ws5.Range("A7").Value = "'4379"
$ws5.Range("B7").Value = "'0"
$ws5.Range("C7").Value = "'10.00%"
$ws5.Range("A8").Value = "'4444"
$ws5.Range("B8").Value = "'0"
$ws5.Range("C8").Value = "'20.00%"
$ws5.Range("A9").Value = "'4446"
$ws5.Range("B9").Value = "'1"
$ws5.Range("C9").Value = "'10.00%"
$ws5.Range("A10").Value = "'4448"
$ws5.Range("B10").Value = "'0"
$ws5.Range("C10").Value = "'40.00%"
$ws5.Range("A11").Value = "'4525"
$ws5.Range("B11").Value = "'"
$ws5.Range("C11").Value = "'"
$ws5.Range("A12").Value = "'4528"
$ws5.Range("B12").Value = "'0"
$ws5.Range("C12").Value = "'10.00%"
$ws5.Range("A13").Value = "'4530"
$ws5.Range("B13").Value = "'0"
$ws5.Range("C13").Value = "'20.00%"
$ws5.Range("A14").Value = "'4537"
$ws5.Range("B14").Value = "'1"
$ws5.Range("C14").Value = "'10.00%"
$ws5.Range("A15").Value = "'4538"
$ws5.Range("B15").Value = "'"
$ws5.Range("C15").Value = "'"
$ws5.Range("A16").Value = "'4539"
$ws5.Range("B16").Value = "'0"
$ws5.Range("C16").Value = "'30.00%"
$ws5.Range("A17").Value = "'4582"
$ws5.Range("B17").Value = "'0"
$ws5.Range("C17").Value = "'20.00%"
$ws5.Range("A18").Value = "'4585"
$ws5.Range("B18").Value = "'0"
$ws5.Range("C18").Value = "'20.00%"
$ws5.Range("A19").Value = "'4588"
$ws5.Range("B19").Value = "'0"
$ws5.Range("C19").Value = "'30.00%"
$ws5.Range("A20").Value = "'4671"
$ws5.Range("B20").Value = "'"
$ws5.Range("C20").Value = "'"
$ws5.Range("A21").Value = "'4675"
$ws5.Range("B21").Value = "'"
$ws5.Range("C21").Value = "'"

# ---------------------------------------------------------------------------
# 2. Clean up "ODI Batting Extra": drop cells that only ever held an empty
#    placeholder string (no real scraped value for that match/column)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ODI Batting Extra")

$ws4.Range("E4").ClearContents()
$ws4.Range("E6").ClearContents()
$ws4.Range("C8").ClearContents()
$ws4.Range("D8").ClearContents()
$ws4.Range("E8").ClearContents()
$ws4.Range("B10").ClearContents()
$ws4.Range("C10").ClearContents()
$ws4.Range("D10").ClearContents()
$ws4.Range("E10").ClearContents()
$ws4.Range("E13").ClearContents()
$ws4.Range("B14").ClearContents()
$ws4.Range("C14").ClearContents()
$ws4.Range("D14").ClearContents()
$ws4.Range("E14").ClearContents()
$ws4.Range("C15").ClearContents()
$ws4.Range("D15").ClearContents()
$ws4.Range("E15").ClearContents()
$ws4.Range("C17").ClearContents()
$ws4.Range("D17").ClearContents()
$ws4.Range("E17").ClearContents()
$ws4.Range("B19").ClearContents()
$ws4.Range("C19").ClearContents()
$ws4.Range("D19").ClearContents()
$ws4.Range("E19").ClearContents()
$ws4.Range("B21").ClearContents()
$ws4.Range("C21").ClearContents()
$ws4.Range("D21").ClearContents()
$ws4.Range("E21").ClearContents()

Write-Host "Added ODI Bowling Extra sheet and cleaned ODI Batting Extra placeholders."
